$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1!A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value()
$text = $text.Replace("1000 Bs = 1.8 = 6547.94 pesos", "1000 Bs = 1.81 = 6605.65 pesos")
$text = $text.Replace("6547.94 pesos = 1.79 = 940.55 Bs", "6605.65 pesos = 1.81 = 968.41 Bs")
$cell.Value = $text

# --- Update the rate figures on the "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 551.8
$wsTasas.Range("O10").Value = 3645
$wsTasas.Range("N12").Value = 3650
$wsTasas.Range("O12").Value = 535.1
